# Adds two new columns (I: "I0", J: "IF") to the sheet, mirroring the
# existing header style from H1 and filling in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells: copy H1's formatting (bold, bordered, centered) onto
# --- I1/J1 so the new headers match the look of the existing ones, then
# --- set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-63: (row, I value, J value) ---
$data = @(
    @(2,8,8), @(3,8,8), @(4,1,1), @(5,7,8), @(6,8,8), @(7,2,3), @(8,8,8), @(9,9,9), @(10,5,6),
    @(11,6,6), @(12,1,1), @(13,6,6), @(14,8,8), @(15,7,8), @(16,1,1), @(17,1,1), @(18,7,7), @(19,8,8),
    @(20,1,1), @(21,10,10), @(22,1,2), @(23,3,3), @(24,7,8), @(25,1,1), @(26,6,6), @(27,10,10), @(28,8,8),
    @(29,4,5), @(30,9,9), @(31,8,8), @(32,9,9), @(33,8,8), @(34,6,7), @(35,1,1), @(36,8,8), @(37,9,9),
    @(38,6,6), @(39,7,7), @(40,6,7), @(41,6,6), @(42,8,8), @(43,8,8), @(44,6,6), @(45,7,8), @(46,8,8),
    @(47,8,8), @(48,8,8), @(49,7,8), @(50,7,8), @(51,9,9), @(52,4,5), @(53,6,6), @(54,7,7), @(55,7,7),
    @(56,8,8), @(57,6,6), @(58,5,6), @(59,4,4), @(60,6,6), @(61,5,5), @(62,4,4), @(63,3,3)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
